$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "JSU(-1.0739450830424986, 1.3092572159085827, 0.09492659788669949, 3.37925668276948)"
$ws.Range("C2").Value = "JSU(-1.2615397936371018, 1.001447542101232, 2.73209754060362, 3.6575022649221447)"
$ws.Range("D2").Value = "NIG(0.8934692927876992, 0.6194881118464688, 1.119879580525911, 2.938698980949743)"
$ws.Range("E2").Value = "NIG(0.8620945837649352, 0.6019800251659057, 5.156222081104403, 5.4871546565516836)"
